$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("D").Insert()
Write-Host "After insert, UsedRange:" $ws.UsedRange.Address()
Write-Host "D7:" ($ws.Range("D7").Value)
Write-Host "E7:" ($ws.Range("E7").Value)
Write-Host "D7 string:" "$($ws.Range('D7').Value)"
Write-Host "E7 string:" "$($ws.Range('E7').Value)"
